$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 16 (shifts old rows 16-26 down to 17-27),
# pushing the new "MAXILASE" item into its alphabetically-sorted slot.
$ws.Rows.Item(16).Insert()

# --- Row 16 (new item row: MAXILASE) ---
$ws.Range("A16:B16").Merge()
$ws.Range("C16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()
$ws.Range("N16:O16").Merge()
$ws.Rows.Item(16).RowHeight = 25.5

$ws.Range("A16").Value = 10
$ws.Range("C16").Value = "MAXILASE 200 CEIP UNIT/ML SYRUP 100ML"
$ws.Range("H16").Value = "3:0"
$ws.Range("L16").Value = "1"
$ws.Range("N16").Value = "57.00"
$ws.Range("P16").Value = "57.0000"
$ws.Range("Q16").Value = "1:0"

# --- Renumber the serial column for the shifted-down item rows (17-25) ---
$ws.Range("A17").Value = 11
$ws.Range("A18").Value = 12
$ws.Range("A19").Value = 13
$ws.Range("A20").Value = 14
$ws.Range("A21").Value = 15
$ws.Range("A22").Value = 16
$ws.Range("A23").Value = 17
$ws.Range("A24").Value = 18
$ws.Range("A25").Value = 19

# --- Update the grand-total cell (row 26, was row 25) to include the new item ---
$ws.Range("P26").Value = 796.95000000000005

# --- Update the footer timestamp (row 27, was row 26) ---
$ws.Range("A27").Value = "Wednesday, 13 August, 2025 10:32 AM"
